$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new column C with header date, copying the style/format from B1
$ws.Range("B1").Copy()
$ws.Range("C1").PasteSpecial(-4122) # xlPasteFormats
$ws.Range("C1").Value = "13-01-2023"

# Reorder rows: stock rows move up to rows 2-4, avg to row 5, total to row 6
# Row 2: Alpha Acciones
$ws.Range("A2").Value = "Alpha Acciones"
$ws.Range("B2").Value = 165685.73
$ws.Range("C2").Value = 160107.87

# Row 3: Arpenta acciones
$ws.Range("A3").Value = "Arpenta acciones"
$ws.Range("B3").Value = 74942.10000000001
$ws.Range("C3").Value = 75167.3

# Row 4: Pellegrini Acciones
$ws.Range("A4").Value = "Pellegrini Acciones"
$ws.Range("B4").Value = 49659.63
$ws.Range("C4").Value = 37230.68

# Row 5: avg
$ws.Range("A5").Value = "avg"
$ws.Range("B5").Value = 96762.49000000001
$ws.Range("C5").Value = 90835.28

# Row 6: total
$ws.Range("A6").Value = "total"
$ws.Range("B6").Value = 290287.46
$ws.Range("C6").Value = 272505.85
